$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-22 13:17:14"
$wsZh.Range("H2").Value = "2016-03-22 13:17:43"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-22 13:17:19"
$wsDe.Range("H2").Value = "2016-03-22 13:17:50"
